$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.130.98"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.562.77"
$ws.Range("E3").Value = "  +5.05%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").Value = "3.557.78"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +5.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.602"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "694.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "4.127.65"
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "71.288.43"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "3.559.29"
$ws.Range("E18").Value = "  +5.02%  "
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "585.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "3.661.44"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.143"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "0.0₃0756"
$ws.Range("E41").Value = "  +8.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.344"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
